# Apply the updated cryptos list values (prices / volumes / row shifts) scraped on
# Tue Jun  6 09:06:57 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume columns so values such as "277.50" or
# "0.06670" are stored verbatim (matching the source inline-string cells) instead of
# being auto-coerced into numbers and losing trailing zeros / thousands dots.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.727.63'
$ws.Range("E2").Value = '  -4.10%  '

$ws.Range("D3").Value = '1.816.98'
$ws.Range("E3").Value = '  -3.06%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '277.50'
$ws.Range("E5").Value = '  -7.94%  '

$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").Value = '0.5091'
$ws.Range("E7").Value = '  -5.13%  '

$ws.Range("D8").Value = '0.3520'
$ws.Range("E8").Value = '  -6.16%  '

$ws.Range("D9").Value = '44.39'
$ws.Range("E9").Value = '  -2.35%  '

$ws.Range("D10").Value = '0.06670'
$ws.Range("E10").Value = '  -7.25%  '

$ws.Range("D11").Value = '20.07'
$ws.Range("E11").Value = '  -7.10%  '

$ws.Range("D12").Value = '0.8257'
$ws.Range("E12").Value = '  -7.26%  '

$ws.Range("D13").Value = '0.07900'
$ws.Range("E13").Value = '  -3.23%  '

$ws.Range("D14").Value = '1.829.11'
$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").Value = '5.071'
$ws.Range("E15").Value = '  -4.72%  '

$ws.Range("D16").Value = '87.63'
$ws.Range("E16").Value = '  -6.13%  '

$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").Value = '14.07'
$ws.Range("E18").Value = '  -5.35%  '

$ws.Range("D19").Value = '0.000008033'
$ws.Range("E19").Value = '  -5.91%  '

$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").Value = '25.769.15'
$ws.Range("E21").Value = '  -4.07%  '

$ws.Range("D22").Value = '4.740'
$ws.Range("E22").Value = '  -4.95%  '

$ws.Range("E23").Value = '  -5.94%  '

$ws.Range("D24").Value = '6.103'

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = '2.211'
$ws.Range("E25").Value = '  -3.98%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '141.87'
$ws.Range("E26").Value = '  -3.03%  '

$ws.Range("E27").Value = '  -3.10%  '

$ws.Range("D28").Value = '17.11'
$ws.Range("E28").Value = '  -5.47%  '

$ws.Range("D29").Value = '109.31'
$ws.Range("E29").Value = '  -4.11%  '

$ws.Range("D30").Value = '4.329'
$ws.Range("E30").Value = '  -8.37%  '

$ws.Range("D31").Value = '4.231'
$ws.Range("E31").Value = '  -8.36%  '

$ws.Range("D32").Value = '0.08769'
$ws.Range("E32").Value = '  -4.18%  '

$ws.Range("D33").Value = '0.04888'
$ws.Range("E33").Value = '  -2.67%  '

$ws.Range("D34").Value = '0.7265'
$ws.Range("E34").Value = '  -10.64%  '

$ws.Range("D35").Value = '1.135'
$ws.Range("E35").Value = '  -3.53%  '

$ws.Range("D36").Value = '2.869'
$ws.Range("E36").Value = '  -2.73%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '3.127'
$ws.Range("E38").Value = '  -2.55%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.373'
$ws.Range("E39").Value = '  -9.61%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01853'
$ws.Range("E40").Value = '  -5.23%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5179'
$ws.Range("E41").Value = '  -13.99%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.9662'
$ws.Range("E42").Value = '  -9.69%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.225'
$ws.Range("E43").Value = '  -6.16%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '110.92'
$ws.Range("E44").Value = '  -3.72%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '8.009'
$ws.Range("E45").Value = '  -9.94%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4566'
$ws.Range("E47").Value = '  -10.45%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1364'
$ws.Range("E48").Value = '  -8.71%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '36.45'
$ws.Range("E49").Value = '  -3.38%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.199'
$ws.Range("E50").Value = '  -8.04%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.501'
$ws.Range("E51").Value = '  -8.24%  '
